$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = 333.48
$ws.Range("E2").Value = 50.6
$ws.Range("F2").Value = 1.14
$ws.Range("K2").Value = 60.5
$ws.Range("N2").Value = 53.62998959737769

# Row 3
$ws.Range("D3").Value = 78.08
$ws.Range("E3").Value = 43.3
$ws.Range("F3").Value = 2.3
$ws.Range("K3").Value = 56.5
$ws.Range("N3").Value = 53.62998959737769

# Row 4
$ws.Range("D4").Value = 77.78
$ws.Range("E4").Value = 46.3
$ws.Range("F4").Value = 2.28
$ws.Range("K4").Value = 50.1
$ws.Range("N4").Value = 53.62998959737769

# Row 5
$ws.Range("D5").Value = 110.92
$ws.Range("E5").Value = 65.3
$ws.Range("F5").Value = 2.73
$ws.Range("H5").Value = 46
$ws.Range("I5").Value = 40
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 50.1
$ws.Range("N5").Value = 53.62998959737769
